$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.897.02"
$ws.Range("E2").Value = "  -3.57%  "
$ws.Range("D3").Value = "1.855.03"
$ws.Range("E3").Value = "  -2.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4351"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3676"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07476"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9363"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.31%  "
$ws.Range("D12").Value = "1.902.43"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.693"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.419"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06882"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.61%  "
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.70%  "
$ws.Range("D21").Value = "27.878.17"
$ws.Range("E21").Value = "  -3.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.095"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "2.103.39"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.005"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.361"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.731"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08978"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7994"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.818"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.017"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.30%  "
$ws.Range("E35").Value = "  -6.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.004"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.112"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05416"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01953"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.902"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5226"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.989"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1678"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.740"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06702"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4868"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.919"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.08%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.003"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.674"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.96%  "
